$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '46.896.41'
$ws.Range('E2').Value = '  +5.01%  '
$ws.Range('D3').Value = '2.335.07'
$ws.Range('E3').Value = '  +4.05%  '
$ws.Range('E4').Value = '  -0.73%  '
$ws.Range('D5').Value = '305.89'
$ws.Range('E5').Value = '  -0.24%  '
$ws.Range('D6').Value = '97.26'
$ws.Range('E6').Value = '  +3.13%  '
$ws.Range('D7').Value = '0.576'
$ws.Range('E7').Value = '  +1.19%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.57%  '
$ws.Range('E9').Value = '  +3.59%  '
$ws.Range('D10').Value = '35.66'
$ws.Range('E10').Value = '  +2.50%  '
$ws.Range('D11').Value = '0.0808'
$ws.Range('E11').Value = '  +0.97%  '
$ws.Range('D12').Value = '7.40'
$ws.Range('E12').Value = '  +2.88%  '
$ws.Range('E13').Value = '  -0.32%  '
$ws.Range('D14').Value = '2.692.54'
$ws.Range('E14').Value = '  +4.09%  '
$ws.Range('D15').Value = '2.337.04'
$ws.Range('E15').Value = '  +4.18%  '
$ws.Range('D16').Value = '14.12'
$ws.Range('E16').Value = '  +4.13%  '
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('D18').Value = '46.749.46'
$ws.Range('E18').Value = '  +5.24%  '
$ws.Range('D19').Value = '13.67'
$ws.Range('E19').Value = '  +16.01%  '
$ws.Range('D20').Value = '0.0₃0948'
$ws.Range('E20').Value = '  +1.39%  '
$ws.Range('E21').Value = '  +0.00%  '
$ws.Range('D22').Value = '67.20'
$ws.Range('E22').Value = '  +2.90%  '
$ws.Range('D23').Value = '246.35'
$ws.Range('E23').Value = '  +3.61%  '
$ws.Range('D24').Value = '2.98'
$ws.Range('E24').Value = '  +1.38%  '
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('D27').Value = '42.02'
$ws.Range('E27').Value = '  +13.82%  '
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('D29').Value = '9.87'
$ws.Range('E29').Value = '  +1.02%  '
$ws.Range('D30').Value = '20.11'
$ws.Range('E30').Value = '  +0.58%  '
$ws.Range('D31').Value = '5.75'
$ws.Range('E31').Value = '  -1.75%  '
$ws.Range('D32').Value = '152.68'
$ws.Range('E32').Value = '  +2.78%  '
$ws.Range('D33').Value = '0.0818'
$ws.Range('E33').Value = '  +4.61%  '
$ws.Range('D34').Value = '2.59'
$ws.Range('E34').Value = '  -1.11%  '
$ws.Range('D35').Value = '3.17'
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').Value = '0.112'
$ws.Range('E36').Value = '  +2.82%  '
$ws.Range('D37').Value = '0.119'
$ws.Range('E37').Value = '  +1.05%  '
$ws.Range('E38').Value = '  -2.44%  '
$ws.Range('D39').Value = '4.01'
$ws.Range('E39').Value = '  +6.31%  '
$ws.Range('D40').Value = '0.0315'
$ws.Range('E40').Value = '  +5.71%  '
$ws.Range('D41').Value = '3.39'
$ws.Range('E41').Value = '  +1.49%  '
$ws.Range('D42').Value = '13.78'
$ws.Range('E42').Value = '  -8.95%  '
$ws.Range('D43').Value = '0.999'
$ws.Range('E43').Value = '  -0.77%  '
$ws.Range('D44').Value = '1.98'
$ws.Range('E44').Value = '  +11.03%  '
$ws.Range('D45').Value = '1.836.77'
$ws.Range('E45').Value = '  +1.48%  '
$ws.Range('E46').Value = '  +5.17%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').Value = '82.10'
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').Value = '74.25'
$ws.Range('E48').Value = '  +7.39%  '
$ws.Range('D49').Value = '4.95'
$ws.Range('E49').Value = '  +2.86%  '
$ws.Range('D50').Value = '98.66'
$ws.Range('E50').Value = '  +0.76%  '
$ws.Range('D51').Value = '55.05'
$ws.Range('E51').Value = '  +2.05%  '
